$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header font: add "Arial" as the font name (keeps existing bold + white color) ---
$ws.Range("A1:G1").Font.Name = "Arial"

# --- 2. Recolor: new custom colors "TOPAZE" (gold/brown family) and "AMETHYSTE" ---
#     header band   : 4472C4 (blue)       -> 996515
#     "blue" rows    : BDD7EE (light blue) -> FFE5B4
#     "green" rows   : E2EFDA (light green)-> FFD700
# (multi-area ranges only paint their first Area via a single assignment in
#  this host, so walk .Areas explicitly to cover every row group)
$ws.Range("A1:G1").Interior.Color = 1402265

$blueRows = $ws.Range("A2:G3,A5:G7,A9:G11,A13:G15,A17:G19,A21:G21,A23:G23,A25:G29")
foreach ($area in $blueRows.Areas) {
  $area.Interior.Color = 11855359
}

$greenRows = $ws.Range("A4:G4,A8:G8,A12:G12,A16:G16,A20:G20,A22:G22,A24:G24")
foreach ($area in $greenRows.Areas) {
  $area.Interior.Color = 55295
}

# --- 3. Column F (duration / Durée) gets a bit narrower: raw width 10 -> 9 ---
$ws.Columns(6).ColumnWidth = 8.17

# --- 4. Header row: translate / re-case the column titles ---
$ws.Range("A1").Value = "Round"
$ws.Range("B1").Value = "Début"
$ws.Range("C1").Value = "Fin"
$ws.Range("D1").Value = "Équipe 1"
$ws.Range("E1").Value = "Équipe 2"
$ws.Range("F1").Value = "Durée"
$ws.Range("G1").Value = "Phase"
